# Updated cryptos list on Wed Sep 11 05:33:04 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column (D) holds plain-text values (e.g. "508.01"). Excel's
# COM layer auto-coerces a numeric-looking string typed into a General
# cell into a real number, which would flip these cells from text to
# numeric storage. Temporarily mark the column as Text before writing so
# the values round-trip as strings, then restore the original (General)
# style so no stray per-cell formatting is left behind.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

# Row 2 - Bitcoin
$ws.Range("D2").Value = "56.233.96"
$ws.Range("E2").Value = "  -1.22%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "2.323.10"
$ws.Range("E3").Value = "  -1.02%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.09%  "

# Row 5 - BNB
$ws.Range("D5").Value = "508.01"
$ws.Range("E5").Value = "  -2.13%  "

# Row 6 - Solana
$ws.Range("D6").Value = "131.68"
$ws.Range("E6").Value = "  -1.83%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  +0.06%  "

# Row 8 - XRP
$ws.Range("E8").Value = "  -0.71%  "

# Row 9 - Dogecoin
$ws.Range("D9").Value = "0.0994"
$ws.Range("E9").Value = "  -3.86%  "

# Row 10 - TRON
$ws.Range("E10").Value = "  -1.07%  "

# Row 11 - Toncoin
$ws.Range("D11").Value = "5.22"
$ws.Range("E11").Value = "  -0.83%  "

# Row 12 - Cardano
$ws.Range("D12").Value = "0.335"
$ws.Range("E12").Value = "  -1.59%  "

# Row 13 - WrappedliquidstakedEther2.0
$ws.Range("D13").Value = "2.739.79"
$ws.Range("E13").Value = "  -0.74%  "

# Row 14 - Avalanche
$ws.Range("D14").Value = "23.34"
$ws.Range("E14").Value = "  -1.74%  "

# Row 15 - WrappedBTC
$ws.Range("D15").Value = "56.232.68"
$ws.Range("E15").Value = "  -1.06%  "

# Row 16 - ShibaInu
$ws.Range("E16").Value = "  -2.12%  "

# Row 17 - WrappedEther
$ws.Range("D17").Value = "2.322.03"
$ws.Range("E17").Value = "  -0.78%  "

# Row 18 - Chainlink
$ws.Range("D18").Value = "10.33"
$ws.Range("E18").Value = "  -1.06%  "

# Row 19 - BitcoinCash
$ws.Range("D19").Value = "320.85"
$ws.Range("E19").Value = "  -0.14%  "

# Row 20 - Polkadot
$ws.Range("D20").Value = "4.12"
$ws.Range("E20").Value = "  -3.34%  "

# Row 21 - Uniswap
$ws.Range("E21").Value = "  -0.28%  "

# Row 22 - Dai
$ws.Range("D22").Value = "0.999"
$ws.Range("E22").Value = "  -0.28%  "

# Row 23 - Litecoin
$ws.Range("D23").Value = "61.07"
$ws.Range("E23").Value = "  +0.47%  "

# Row 24 - InternetComputer(DFINITY)
$ws.Range("D24").Value = "8.49"
$ws.Range("E24").Value = "  +9.45%  "

# Row 25 - Kaspa
$ws.Range("E25").Value = "  +0.69%  "

# Row 26 - Binance-PegBSC-USD
$ws.Range("D26").Value = "1.00"
$ws.Range("E26").Value = "  +0.02%  "

# Row 27 - Fetch.AI
$ws.Range("E27").Value = "  +3.22%  "

# Row 28 - Monero
$ws.Range("D28").Value = "167.16"
$ws.Range("E28").Value = "  -1.77%  "

# Row 29 - PancakeSwap
$ws.Range("D29").Value = "1.65"
$ws.Range("E29").Value = "  -1.82%  "

# Row 30 - PEPE
$ws.Range("D30").Value = "0.0₃0711"
$ws.Range("E30").Value = "  -4.01%  "

# Row 31 - Aptos
$ws.Range("E31").Value = "  -2.78%  "

# Row 32 - EthereumClassic
$ws.Range("D32").Value = "18.21"
$ws.Range("E32").Value = "  -0.38%  "

# Row 33 - USDe
$ws.Range("E33").Value = "  -0.04%  "

# Row 34 - FirstDigitalUSD
$ws.Range("D34").Value = "1.00"
$ws.Range("E34").Value = "  +0.07%  "

# Row 35 - ImmutableX
$ws.Range("E35").Value = "  -0.05%  "

# Row 36 - NEARProtocol
$ws.Range("D36").Value = "3.90"
$ws.Range("E36").Value = "  -2.28%  "

# Row 37 - SuiNetwork
$ws.Range("D37").Value = "0.876"
$ws.Range("E37").Value = "  -6.38%  "

# Row 38 - OKB
$ws.Range("D38").Value = "38.43"
$ws.Range("E38").Value = "  +2.37%  "

# Row 39 - was Stacks, now Aave (rows 39/40 swapped)
$ws.Range("B39").Value = "Aave"
$ws.Range("C39").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D39").Value = "149.26"
$ws.Range("E39").Value = "  +8.34%  "

# Row 40 - was Aave, now Stacks
$ws.Range("B40").Value = "Stacks"
$ws.Range("C40").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D40").Value = "1.52"
$ws.Range("E40").Value = "  -0.41%  "

# Row 41 - PolygonEcosystemToken
$ws.Range("E41").Value = "  -2.10%  "

# Row 42 - Filecoin
$ws.Range("D42").Value = "3.53"
$ws.Range("E42").Value = "  -1.65%  "

# Row 43 - Bittensor
$ws.Range("D43").Value = "276.06"
$ws.Range("E43").Value = "  -0.21%  "

# Row 44 - RenderToken
$ws.Range("D44").Value = "4.97"
$ws.Range("E44").Value = "  -2.18%  "

# Row 45 - Stellar
$ws.Range("D45").Value = "0.0920"
$ws.Range("E45").Value = "  -1.37%  "

# Row 46 - Hedera
$ws.Range("D46").Value = "0.0491"
$ws.Range("E46").Value = "  -2.54%  "

# Row 47 - Mantle
$ws.Range("E47").Value = "  -1.30%  "

# Row 48 - InjectiveProtocol
$ws.Range("D48").Value = "17.80"
$ws.Range("E48").Value = "  +3.29%  "

# Row 49 - Polygon
$ws.Range("E49").Value = "  -0.90%  "

# Row 50 - VeChain
$ws.Range("E50").Value = "  -2.12%  "

# Row 51 - EnergySwap
$ws.Range("D51").Value = "16.93"
$ws.Range("E51").Value = "  +0.41%  "

# Restore the Price column's original (default/General) style now that
# every text value has been written, so no extra per-cell styling is
# left over from the temporary Text format above.
$priceRange.Style = "Normal"
